# Loan RBI, Variable Instalments
# On the "Repayment schedule" sheet, insert a new (blank) column before
# column N so the existing "Late" / "heading" / "Outstanding" columns
# shift one place to the right (N->O, O->P, P->Q). Keep the new column's
# width consistent with its neighbour, make "Repayment schedule" the
# active sheet/tab, and leave the selection on K14 (matching the
# author's final cursor position).

$wb = $excel.ActiveWorkbook
$ws = $wb.Worksheets.Item("Repayment schedule")

# Remember column M's width so the freshly inserted column N can match
# it (both end up stored as width "11" in the OOXML, same as a manual
# "Insert Column" followed by a drag-to-resize would produce).
$refWidth = $ws.Columns("M:M").ColumnWidth

# Insert a new blank column at N; everything at/after N shifts right.
$ws.Columns("N:N").Insert()
$ws.Columns("N:N").ColumnWidth = $refWidth

# Make "Repayment schedule" the active sheet/tab and park the selection
# on K14.
$ws.Activate()
[void]$ws.Range("K14").Select()
